$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 753, shifting existing rows 753-816 down to 755-818
$ws.Rows("753:754").Insert()

# Populate the two newly inserted rows with the new price data
# Row 753
$ws.Cells.Item(753, 1).Value = 9
$ws.Cells.Item(753, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(753, 3).Value = 'Metropolitana'
$ws.Cells.Item(753, 4).Value = 45132
$ws.Cells.Item(753, 5).Value = 13
$ws.Cells.Item(753, 6).Value = 100112031
$ws.Cells.Item(753, 7).Value = 'Poroto verde'
$ws.Cells.Item(753, 8).Value = 'Magnum'
$ws.Cells.Item(753, 9).Value = 'Primera'
$ws.Cells.Item(753, 10).Value = 70
$ws.Cells.Item(753, 11).Value = 19000
$ws.Cells.Item(753, 12).Value = 21000
$ws.Cells.Item(753, 13).Value = 20000
$ws.Cells.Item(753, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(753, 15).Value = 'Perú'
$ws.Cells.Item(753, 16).Value = 800
$ws.Cells.Item(753, 17).Value = 25
$ws.Cells.Item(753, 18).Value = 'Hortaliza'

# Row 754
$ws.Cells.Item(754, 1).Value = 9
$ws.Cells.Item(754, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(754, 3).Value = 'Metropolitana'
$ws.Cells.Item(754, 4).Value = 45132
$ws.Cells.Item(754, 5).Value = 13
$ws.Cells.Item(754, 6).Value = 100112031
$ws.Cells.Item(754, 7).Value = 'Poroto verde'
$ws.Cells.Item(754, 8).Value = 'Sin especificar'
$ws.Cells.Item(754, 9).Value = 'Primera'
$ws.Cells.Item(754, 10).Value = 52
$ws.Cells.Item(754, 11).Value = 25000
$ws.Cells.Item(754, 12).Value = 27000
$ws.Cells.Item(754, 13).Value = 26000
$ws.Cells.Item(754, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(754, 15).Value = 'Perú'
$ws.Cells.Item(754, 16).Value = 1040
$ws.Cells.Item(754, 17).Value = 25
$ws.Cells.Item(754, 18).Value = 'Hortaliza'
